# Updates the Price (D) and Volume(1h) (E) columns for the cryptocurrency
# rows (2-51) with freshly scraped values, matching the GitHub Actions
# "Updated cryptos list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "26.191.02"; E = "  -0.81%  " },
    @{ Row = 3; D = "1.681.86"; E = "  -0.49%  " },
    @{ Row = 4; D = "1.005"; E = "  -0.57%  " },
    @{ Row = 5; D = "212.40"; E = "  -3.09%  " },
    @{ Row = 6; D = "0.5294"; E = "  -4.01%  " },
    @{ Row = 7; D = "1.005"; E = "  -0.52%  " },
    @{ Row = 8; D = "0.2695"; E = "  -0.53%  " },
    @{ Row = 9; D = "0.06343"; E = "  -2.27%  " },
    @{ Row = 10; D = "21.44"; E = "  -3.02%  " },
    @{ Row = 11; D = "0.07568"; E = "  -0.49%  " },
    @{ Row = 12; D = "1.686.08"; E = "  -0.29%  " },
    @{ Row = 13; D = "4.527"; E = "  -0.86%  " },
    @{ Row = 14; D = "0.5716"; E = "  -1.87%  " },
    @{ Row = 15; D = "0.000008187"; E = "  -3.27%  " },
    @{ Row = 16; D = "66.40"; E = "  +1.49%  " },
    @{ Row = 17; D = "26.254.76"; E = "  -1.05%  " },
    @{ Row = 18; D = "4.886"; E = "  -1.25%  " },
    @{ Row = 19; D = "1.005"; E = "  -0.52%  " },
    @{ Row = 20; D = "10.65"; E = "  -2.96%  " },
    @{ Row = 21; D = "190.77"; E = "  +0.05%  " },
    @{ Row = 22; D = "6.253"; E = "  -0.02%  " },
    @{ Row = 23; D = "1.006"; E = "  -0.51%  " },
    @{ Row = 24; D = "149.05"; E = "  -0.43%  " },
    @{ Row = 25; D = "0.1265"; E = "  -3.78%  " },
    @{ Row = 26; D = "7.682"; E = "  -3.42%  " },
    @{ Row = 27; D = "15.96"; E = "  +0.93%  " },
    @{ Row = 28; D = "0.06527"; E = "  +3.44%  " },
    @{ Row = 29; D = "1.345"; E = "  -4.71%  " },
    @{ Row = 30; D = "1.292"; E = "  -2.76%  " },
    @{ Row = 31; D = "3.561"; E = "  -0.82%  " },
    @{ Row = 32; D = "3.566"; E = "  -0.61%  " },
    @{ Row = 33; D = "1.674"; E = "  -0.11%  " },
    @{ Row = 34; D = "1.014"; E = "  -3.11%  " },
    @{ Row = 35; D = "0.6119"; E = "  -2.05%  " },
    @{ Row = 36; D = "2.415"; E = "  +0.20%  " },
    @{ Row = 37; D = "2.724"; E = "  +0.02%  " },
    @{ Row = 38; D = "6.197"; E = "  -0.77%  " },
    @{ Row = 39; D = "0.01618"; E = "  -1.64%  " },
    @{ Row = 40; D = "1.102.14"; E = "  -1.69%  " },
    @{ Row = 41; D = "0.8722"; E = "  -1.11%  " },
    @{ Row = 42; D = "1.007"; E = "  -0.92%  " },
    @{ Row = 43; D = "100.20"; E = "  -0.70%  " },
    @{ Row = 44; D = "1.834.57"; E = "  -0.36%  " },
    @{ Row = 45; D = "0.00000000110"; E = "  -4.75%  " },
    @{ Row = 46; D = "57.38"; E = "  -0.29%  " },
    @{ Row = 47; D = "1.007"; E = "  +0.06%  " },
    @{ Row = 48; D = "8.038"; E = "  -2.18%  " },
    @{ Row = 49; D = "0.05267"; E = "  -0.30%  " },
    @{ Row = 50; D = "0.4273"; E = "  -0.75%  " },
    @{ Row = 51; D = "5.994"; E = "  -1.74%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    # Column D holds price text that can look like a plain number
    # (e.g. "212.40"). Force text format before assigning so Excel keeps
    # the exact original formatting (trailing zeros, many decimals, etc.)
    # instead of silently coercing it to a numeric value, then clear the
    # formatting override again so the cell keeps its original (default)
    # style.
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
    $dCell.ClearFormats()

    # Column E is always a padded percentage string (never numeric-looking)
    # so it can be assigned directly.
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $u.E
}
